# Fix CSEI93 source codes: correct the recorded answer options on the
# "question_answers" sheet and the resulting expected scores on the
# "outputs" sheet.

$wb = $excel.ActiveWorkbook

$qa = $wb.Worksheets.Item("question_answers")
$out = $wb.Worksheets.Item("outputs")

# Corrected answers (question row -> new answer option) on question_answers!B
$answerFixes = @{
    2  = "2"
    3  = "1"
    5  = "2"
    6  = "2"
    8  = "1"
    9  = "1"
    10 = "1"
    12 = "1"
    13 = "1"
    15 = "2"
    16 = "2"
    19 = "1"
    26 = "2"
    27 = "2"
    28 = "1"
    30 = "2"
    31 = "2"
    32 = "1"
    35 = "1"
    38 = "1"
    39 = "2"
    40 = "2"
    42 = "1"
    43 = "1"
    44 = "1"
    45 = "1"
    46 = "1"
    47 = "1"
    49 = "2"
    50 = "2"
    52 = "2"
    56 = "2"
    58 = "1"
    59 = "1"
}

foreach ($row in $answerFixes.Keys) {
    $cell = $qa.Range("B$row")
    # Keep the cell stored as text (matches the original inline-string type)
    # instead of letting Excel auto-convert the numeric-looking "1"/"2" values.
    $cell.NumberFormat = "@"
    $cell.Value = $answerFixes[$row]
}

# Recomputed expected scores on outputs!B following the corrected answers
$out.Range("B2").Value = 12   # general_self_steem
$out.Range("B3").Value = 2    # family_self_steem
$out.Range("B6").Value = 6    # life_scale
$out.Range("B7").Value = 22   # total
